$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: shift rows 6..18 down to 7..19 (bottom-up to avoid clobbering) ---
# old row 18 -> row 19
$ws.Cells.Item(19,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(19,2).Value = "【急募】BOXファイルをGASでkintoneに自動同期したい"
$ws.Cells.Item(19,3).Value = "システム開発"
$ws.Cells.Item(19,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(19,5).Value = "期限情報なし"
$ws.Cells.Item(19,6).Value = "https://www.lancers.jp/work/detail/5487010"
$ws.Cells.Item(19,7).Value = 10

# old row 17 -> row 18
$ws.Cells.Item(18,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(18,2).Value = "【急募】Klaviyoスパム対策とドメイン解決の専門家募集"
$ws.Cells.Item(18,3).Value = "システム開発"
$ws.Cells.Item(18,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(18,5).Value = "期限情報なし"
$ws.Cells.Item(18,6).Value = "https://www.lancers.jp/work/detail/5486673"
$ws.Cells.Item(18,7).Value = 13

# old row 16 -> row 17
$ws.Cells.Item(17,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(17,2).Value = "【急募】SSLエラー解決のための専門家を探しています"
$ws.Cells.Item(17,3).Value = "システム開発"
$ws.Cells.Item(17,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(17,5).Value = "期限情報なし"
$ws.Cells.Item(17,6).Value = "https://www.lancers.jp/work/detail/5486960"
$ws.Cells.Item(17,7).Value = 13

# old row 15 -> row 16
$ws.Cells.Item(16,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(16,2).Value = "【急募】外部CTOを探しています!"
$ws.Cells.Item(16,3).Value = "システム開発"
$ws.Cells.Item(16,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(16,5).Value = "期限情報なし"
$ws.Cells.Item(16,6).Value = "https://www.lancers.jp/work/detail/5486956"
$ws.Cells.Item(16,7).Value = 18

# old row 14 -> row 15
$ws.Cells.Item(15,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(15,2).Value = "【急募】ECサービスのメール送信障害調査・改善支援"
$ws.Cells.Item(15,3).Value = "システム開発"
$ws.Cells.Item(15,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(15,5).Value = "期限情報なし"
$ws.Cells.Item(15,6).Value = "https://www.lancers.jp/work/detail/5487035"
$ws.Cells.Item(15,7).Value = 18

# old row 13 -> row 14
$ws.Cells.Item(14,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(14,2).Value = "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集"
$ws.Cells.Item(14,3).Value = "システム開発"
$ws.Cells.Item(14,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(14,5).Value = "期限情報なし"
$ws.Cells.Item(14,6).Value = "https://www.lancers.jp/work/detail/5486471"
$ws.Cells.Item(14,7).Value = 25

# old row 12 -> row 13
$ws.Cells.Item(13,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(13,2).Value = "【募集】PHP + MySQLでのcron用スクリプト作成依頼"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5487013"
$ws.Cells.Item(13,7).Value = 50
$ws.Cells.Item(13,8).Value = "◇MySQL ○PHP"

# old row 11 -> row 12
$ws.Cells.Item(12,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(12,2).Value = "初心者向けダンススクールの問い合わせフォームを置き換える/拡張するチャットボット開発"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5486851"
$ws.Cells.Item(12,7).Value = 75
$ws.Cells.Item(12,8).Value = "◆開発"

# old row 10 -> row 11
$ws.Cells.Item(11,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(11,2).Value = "【Java/講師】新入社員研修のサブ講師募集"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5477144"
$ws.Cells.Item(11,7).Value = 78
$ws.Cells.Item(11,8).Value = "★Java"

# old row 9 -> row 10
$ws.Cells.Item(10,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(10,2).Value = "【業務改善】訪問業務に特化したスケジュール/介入実績管理Webシステム構築"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5486583"
$ws.Cells.Item(10,7).Value = 85
$ws.Cells.Item(10,8).Value = "◇業務改善"

# old row 8 -> row 9
$ws.Cells.Item(9,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(9,2).Value = "【急募】ガイドと旅行者をつなぐマッチングサイト開発"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5475198"
$ws.Cells.Item(9,7).Value = 93
$ws.Cells.Item(9,8).Value = "◆開発 ◇サイト"

# old row 7 -> row 8
$ws.Cells.Item(8,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(8,2).Value = "【急募】FlutterでのSNS風アプリ開発をお願いします(Firebase想定)"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5487157"
$ws.Cells.Item(8,7).Value = 93
$ws.Cells.Item(8,8).Value = "◆開発 ◇アプリ"

# old row 6 -> row 7
$ws.Cells.Item(7,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(7,2).Value = "【急募】WordPressサイト再構築+LINE・予約連携+顧客管理機能構築|テーマ指定あり|"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5486856"
$ws.Cells.Item(7,7).Value = 93
$ws.Cells.Item(7,8).Value = "◇サイト ○WordPress"

# --- Step 2: write new row 6 content ---
$ws.Cells.Item(6,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(6,2).Value = "【業務自動化】国際郵便マイページの配送ラベル一括印刷の自動化ツール開発"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5487449"
$ws.Cells.Item(6,7).Value = 205
$ws.Cells.Item(6,8).Value = "◆ツール,開発"

# --- Step 3: update timestamp (column A) for rows 2..5 ---
$ws.Cells.Item(2,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(3,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(4,1).Value = "2026-02-07 01:49:50"
$ws.Cells.Item(5,1).Value = "2026-02-07 01:49:50"

# --- Step 4: apply hyperlink style to new F6 cell (copy from existing hyperlink cell) ---
$ws.Cells.Item(6,6).Style = $ws.Cells.Item(2,6).Style

# --- Step 5: rebuild hyperlinks collection in row order (F2..F19) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5460562") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5487324") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460563") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5486863") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5487449") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5486856") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5487157") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5475198") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5486583") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5477144") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5486851") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5487013") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5486471") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5487035") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5486956") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5486960") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5486673") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5487010") | Out-Null

Write-Host "Row insert + shift complete"